# "save data done + era data updated"
# Add a new "Save" column (H) to the sheet: a header cell matching the
# existing header formatting, plus one numeric flag per data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Give H1 the same direct formatting as the other header cells (bold,
# bordered, centered) by copying G1's format onto it, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial($xlPasteFormats)
$ws.Range("H1").Value = "Save"

# New "Save" flag values for each data row (2-7).
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 0

$excel.CutCopyMode = $false
